# Commit: "Added company website scraping for emails"
# Adds a third "Shop Email" header column (C1) next to the existing
# "Shop Name" (A1) and "Shop Website" (B1) headers, and moves the
# selection to the newly added header cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Shop Email"

$ws.Range("C1").Select()
